$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($sheet, $cellRef, $val)
    $cell = $sheet.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws 'D2' '28.156.50'
$ws.Range('E2').Value = '  -1.32%  '
Set-TextValue $ws 'D3' '1.818.15'
$ws.Range('E3').Value = '  +1.12%  '
Set-TextValue $ws 'D4' '1.006'
$ws.Range('E4').Value = '  +0.44%  '
Set-TextValue $ws 'D5' '329.53'
$ws.Range('E5').Value = '  -0.21%  '
Set-TextValue $ws 'D6' '1.003'
$ws.Range('E6').Value = '  +0.25%  '
Set-TextValue $ws 'D7' '0.4419'
$ws.Range('E7').Value = '  -0.51%  '
Set-TextValue $ws 'D8' '0.3737'
$ws.Range('E8').Value = '  -1.19%  '
Set-TextValue $ws 'D9' '44.72'
$ws.Range('E9').Value = '  -1.99%  '
Set-TextValue $ws 'D10' '0.07720'
$ws.Range('E10').Value = '  +0.62%  '
Set-TextValue $ws 'D11' '1.117'
$ws.Range('E11').Value = '  -2.97%  '
Set-TextValue $ws 'D12' '1.004'
$ws.Range('E12').Value = '  +0.09%  '
Set-TextValue $ws 'D13' '21.99'
$ws.Range('E13').Value = '  -3.49%  '
Set-TextValue $ws 'D14' '6.306'
$ws.Range('E14').Value = '  -0.21%  '
Set-TextValue $ws 'D15' '7.524'
$ws.Range('E15').Value = '  +0.51%  '
Set-TextValue $ws 'D16' '1.830.01'
$ws.Range('E16').Value = '  +2.69%  '
Set-TextValue $ws 'D17' '93.06'
$ws.Range('E17').Value = '  +10.87%  '
Set-TextValue $ws 'D18' '0.00001082'
$ws.Range('E18').Value = '  -1.11%  '
Set-TextValue $ws 'D19' '0.06476'
$ws.Range('E19').Value = '  -3.62%  '
Set-TextValue $ws 'D20' '1.003'
$ws.Range('E20').Value = '  -0.29%  '
Set-TextValue $ws 'D21' '17.56'
$ws.Range('E21').Value = '  -0.41%  '
Set-TextValue $ws 'D22' '6.301'
$ws.Range('E22').Value = '  -0.09%  '
Set-TextValue $ws 'D23' '0.5392'
$ws.Range('E23').Value = '  +0.42%  '
Set-TextValue $ws 'D24' '28.256.74'
$ws.Range('E24').Value = '  -0.94%  '
Set-TextValue $ws 'D25' '11.68'
$ws.Range('E25').Value = '  -0.15%  '
Set-TextValue $ws 'D26' '2.114'
$ws.Range('E26').Value = '  -13.11%  '
Set-TextValue $ws 'D27' '20.66'
$ws.Range('E27').Value = '  -0.89%  '
Set-TextValue $ws 'D28' '155.23'
$ws.Range('E28').Value = '  +1.30%  '
$ws.Range('B29').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C29').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws 'D29' '2.035.59'
$ws.Range('E29').Value = '  +2.37%  '
$ws.Range('B30').Value = 'LidoDAOToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws 'D30' '2.339'
$ws.Range('E30').Value = '  -3.59%  '
Set-TextValue $ws 'D31' '127.93'
$ws.Range('E31').Value = '  -2.76%  '
Set-TextValue $ws 'D32' '1.193'
$ws.Range('E32').Value = '  -10.87%  '
Set-TextValue $ws 'D33' '5.846'
$ws.Range('E33').Value = '  -1.69%  '
Set-TextValue $ws 'D34' '0.09267'
$ws.Range('E34').Value = '  -0.66%  '
Set-TextValue $ws 'D35' '3.667'
$ws.Range('E35').Value = '  -7.78%  '
Set-TextValue $ws 'D36' '12.96'
$ws.Range('E36').Value = '  +4.68%  '
Set-TextValue $ws 'D37' '0.02331'
$ws.Range('E37').Value = '  -0.91%  '
Set-TextValue $ws 'D38' '0.2175'
$ws.Range('E38').Value = '  -3.98%  '
Set-TextValue $ws 'D39' '5.156'
$ws.Range('E39').Value = '  -2.43%  '
Set-TextValue $ws 'D40' '0.6563'
$ws.Range('E40').Value = '  -3.13%  '
Set-TextValue $ws 'D41' '0.06163'
$ws.Range('E41').Value = '  -2.97%  '
Set-TextValue $ws 'D42' '1.204'
$ws.Range('E42').Value = '  -0.68%  '
Set-TextValue $ws 'D43' '8.112'
$ws.Range('E43').Value = '  -0.79%  '
$ws.Range('E44').Value = '  -0.38%  '
Set-TextValue $ws 'D45' '13.94'
$ws.Range('E45').Value = '  -0.90%  '
Set-TextValue $ws 'D46' '1.389'
$ws.Range('E46').Value = '  -4.28%  '
Set-TextValue $ws 'D47' '0.6061'
$ws.Range('E47').Value = '  -2.01%  '
Set-TextValue $ws 'D48' '3.768'
$ws.Range('E48').Value = '  -1.53%  '
Set-TextValue $ws 'D49' '2.041'
$ws.Range('E49').Value = '  -0.34%  '
Set-TextValue $ws 'D50' '126.62'
$ws.Range('E50').Value = '  -1.32%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws 'D51' '0.06988'
$ws.Range('E51').Value = '  -0.19%  '
